# correção nos dados e inicio da analise PNAD 2009
#
# The rows that only contained a category-header label (no numeric data),
# plus the two footnote-only rows at the bottom, are removed. Excel
# shifts every following row up when a row is deleted, which reproduces
# the target layout (and — since the shared-string table is rebuilt by
# Excel — also naturally drops the now-unused header/footnote strings
# from xl/sharedStrings.xml).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete bottom-most rows first so earlier row numbers stay valid.
$ws.Rows.Item(35).Delete()   # "(1) inclusive as pessoas de cor ou raça ..."
$ws.Rows.Item(34).Delete()   # "fonte: ibge, diretoria de pesquisas, ..."
$ws.Rows.Item(27).Delete()   # "classes de rendimento mensal domiciliar per capita"
$ws.Rows.Item(19).Delete()   # "nível de instrução"
$ws.Rows.Item(13).Delete()   # "grupos de idade"
$ws.Rows.Item(8).Delete()    # "cor ou raça"
$ws.Rows.Item(5).Delete()    # "sexo"
